$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of an existing header cell (e.g. AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the Wins/Losses/Ties values for each data row (2 through 43)
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 84  # AD
    $ws.Cells.Item($r, 31).Value = 78  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
